$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph, built by copy/pasting the existing bold
#    "Play 6 Wild Sharks..." paragraph from the end of the document so
#    that it keeps the same run layout (leading empty run + bold run),
#    then re-texting the runs.
# ------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$boldHeaderPara = $d.Paragraphs($lastParaIndex - 1)
$boldHeaderPara.Range.Copy()

$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"
$metaPara.Range.Paste()

# Re-text the (still bold) run that now spans the whole pasted paragraph
# so it just reads "Meta description".
$metaParaAgain = $d.Paragraphs(2)
$boldRange = $d.Range($metaParaAgain.Range.Start, $metaParaAgain.Range.End)
$boldRange.Text = "Meta description"

# Append the plain (non-bold) remainder right before the paragraph mark.
$metaParaFinal = $d.Paragraphs(2)
$insertPos = $metaParaFinal.Range.End - 1
$tailRange = $d.Range($insertPos, $insertPos)
$tailRange.InsertAfter(": Experience 6 Wild Sharks, the first slot game with over 14,000 unique Wild Choice configurations. Play now for free and enjoy immersive graphics and impressive RTP.")
$tailRange.Font.Bold = 0

# ------------------------------------------------------------------
# 2) Drop the trailing duplicate "Play 6 Wild Sharks..." bold heading
#    paragraph near the end of the document.
# ------------------------------------------------------------------
$dupIndex = $d.Paragraphs.Count - 1
$dupPara = $d.Paragraphs($dupIndex)
$dupPara.Range.Delete()

# ------------------------------------------------------------------
# 3) Replace the final (italic) paragraph's text with the new image
#    generation prompt, preserving its italic formatting.
# ------------------------------------------------------------------
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End - 1)
$finalRange.Text = 'Prompt: Create a cartoon-style feature image for the game "6 Wild Sharks" featuring a happy Maya warrior with glasses. The image should convey the thrill and excitement of hunting sharks while highlighting the game''s unique Wild Choice mechanics and the customizable wild symbols that players can use to increase their chances of winning big. The Maya warrior should be seen holding a fishing rod and standing on a boat with 6 wild sharks jumping out of the water in the background. The image should be bright and colorful, with bold text reading "6 Wild Sharks" and "Customize Your Game" to showcase the game''s innovative features.'
